$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# Update the "Run" (yes/no) values in column D
$ws.Range("D2").Value = "no"
$ws.Range("D3").Value = "yes"
$ws.Range("D6").Value = "yes"

# Update the active selection on the sheet to D4
$ws.Activate()
$ws.Range("D4").Select()
